# g8.3 - correção nas colunas
#
# The sheet originally stored two stacked blocks of the same products:
#   rows 2-10  : "Categoria" = 2025-2024, "Valor" = the 2025-2024 delta
#   rows 11-19 : "Categoria" = 2025/1997, "Valor" = the 2025/1997 delta
#
# The corrected layout merges both blocks side-by-side into a single
# 10-row table: column B now holds the 2025-2024 delta and column C
# now holds the 2025/1997 delta, and the two header labels are renamed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the existing values before anything is changed/shifted ---
$catDelta = @{}
$totalDelta = @{}
for ($r = 2; $r -le 10; $r++) {
    $catDelta[$r] = $ws.Cells.Item($r, 3).Value()
}
for ($r = 11; $r -le 19; $r++) {
    $totalDelta[$r - 9] = $ws.Cells.Item($r, 3).Value()
}

# --- rename the header labels ---
$ws.Range("B1").Value = "atual-ano anterior"
$ws.Range("C1").Value = "atual/1997"

# --- rewrite column B (2025-2024 delta) and column C (2025/1997 delta) ---
for ($r = 2; $r -le 10; $r++) {
    $bVal = $catDelta[$r]
    if ($null -eq $bVal -or $bVal -eq "") {
        $ws.Cells.Item($r, 2).ClearContents()
    } else {
        $ws.Cells.Item($r, 2).Value = $bVal
    }
    $ws.Cells.Item($r, 3).Value = $totalDelta[$r]
}

# --- the second (now redundant) block of rows is no longer needed ---
$ws.Rows("11:19").Delete()
